$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BDSBaPCF")

# The existing "hydrogen" plant type is being split into two distinct
# plant types: rename the current row to "hydrogen combustion turbine"
# and add a new row underneath for "hydrogen combined cycle".
$ws.Range("A24").Value = "hydrogen combustion turbine"

$ws.Range("A25").Value = "hydrogen combined cycle"
$ws.Range("B25").Value = 0

# A19:A23 no longer carry the special formatting that used to be shared
# with the "hydrogen" row - only the hydrogen-related rows (now A24:A25)
# keep it, with a slightly different look (pure black font + vertically
# centered text).
$ws.Range("A19:A23").ClearFormats()

$ws.Range("A24").Font.Color = 0
$ws.Range("A24").VerticalAlignment = -4108

# Copy the finished A24 formatting (and the B24 number format/fill) down
# onto the new row so both hydrogen rows look identical.
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B24").Copy()
$ws.Range("B25").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# Make BDSBaPCF the active sheet/tab, with D32 as the last-used selection,
# matching where the editor ended up after making the change.
$ws.Activate()
$ws.Range("D32").Select()
